# Continuing update on prep columns: add a new "s1cDNAProtocol" column (G)
# with value "E7420L" for every data row, and select the new data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column G
$ws.Range("G1").Value = "s1cDNAProtocol"

# Fill values for rows 2-21
$ws.Range("G2:G21").Value = "E7420L"
$ws.Range("G2:G21").WrapText = $true

# Update the selection to mirror the authored change
$ws.Range("G2:G21").Select()
